# Generate Report for Handoff
# Update status/handoff-timestamp for the "ffbfa69e-..." file now that it
# has been handed off, across the Overview summary sheet and each
# per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the ffbfa69e-... file -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("C3").Value = "Ready for handoff"   # de-de status
$wsOverview.Range("D3").Value = "2016-12-11 20:12:55" # Latest Handoff Date

# --- zh-cn sheet: row 3 is the ffbfa69e-... file ---------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"    # Status
$wsZhCn.Range("E3").Value = "2016-03-11 20:12:52"  # Latest Handoff Datetime

# --- de-de sheet: row 3 is the ffbfa69e-... file ---------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"    # Status
$wsDeDe.Range("E3").Value = "2016-03-11 20:12:55"  # Latest Handoff Datetime
